$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "38.397.08"
$ws.Range("E2").Value = "  +1.88%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.077.59"
$ws.Range("E3").Value = "  +1.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.38%  "

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.53%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.386.61"
$ws.Range("E12").Value = "  +2.00%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +2.65%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.39%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +0.95%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.068.86"
$ws.Range("E17").Value = "  +1.58%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "38.311.28"
$ws.Range("E18").Value = "  +1.65%  "

# Row 19 - Litecoin
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.70%  "

# Row 20 - Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.08%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.16%  "

# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.72%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.39%  "

# Row 27 - Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.38%  "

# Row 28 - Kaspa
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.135"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.90%  "

# Row 29 - EthereumClassic
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.73%  "

# Row 30 - ImmutableX
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.38%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -0.62%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  +4.73%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.77%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +2.61%  "

# Row 35 - Hedera
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +1.07%  "

# Row 37 - THORChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.67%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +3.15%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  +0.04%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.82%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.537.58"
$ws.Range("E41").Value = "  +0.63%  "

# Row 42 - Aave
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.87%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +2.30%  "

# Row 44 - Cronos
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.76%  "

# Row 45 - HuobiToken
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0921"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.37%  "

# Row 46 - FraxShare
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.70"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.89%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +0.27%  "

# Row 48 - FTXToken
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.69%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +1.92%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +1.52%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.275.67"
$ws.Range("E51").Value = "  +2.14%  "
